$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: grand total (average) of the k column (J) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# --- Summary rows 14-17: labels (new shared strings) + aggregate formulas ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold, slightly larger, vertically centered font for the summary labels' values
$summary = $ws.Range("B14:B17")
$summary.Font.Bold = $true
$summary.VerticalAlignment = -4108
$summary.Font.Size = 12

# Bold the new grand-total cell to match
$ws.Range("J12").Font.Bold = $true

# Rows grew a bit taller to fit the larger font
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# --- Page setup (portrait, A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ends up resting on the new total cell ---
$ws.Range("J12").Select()
